$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.627.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.46%  "

$ws.Range("D3").Value = "'2.678.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.12%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'597.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.81%  "

$ws.Range("D6").Value = "'144.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.92%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'2.677.23"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.11%  "

$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("D11").Value = "'5.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("E14").Value = "  +2.70%  "

$ws.Range("D15").Value = "'3.150.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.90%  "

$ws.Range("D16").Value = "'63.517.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "

$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").Value = "'2.677.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.65%  "

$ws.Range("E19").Value = "  +4.75%  "

$ws.Range("D20").Value = "'4.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.51%  "

$ws.Range("D21").Value = "'340.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'6.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.51%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").Value = "'67.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("D25").Value = "'1.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.42%  "

$ws.Range("D26").Value = "'1.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.51%  "

$ws.Range("E27").Value = "  +1.27%  "

$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("D29").Value = "'540.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +19.15%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +13.14%  "

$ws.Range("D33").Value = "'2.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.93%  "

$ws.Range("D34").Value = "'0.0₃0817"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("D35").Value = "'172.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("D36").Value = "'5.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.32%  "

$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").Value = "'19.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.13%  "

$ws.Range("E40").Value = "  +9.63%  "

$ws.Range("D41").Value = "'174.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.23%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "'40.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("D45").Value = "'22.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.17%  "

$ws.Range("E46").Value = "  +6.27%  "

$ws.Range("D47").Value = "'0.637"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "

$ws.Range("E48").Value = "  +3.18%  "

$ws.Range("D49").Value = "'0.0965"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").Value = "'18.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.34%  "

$ws.Range("E51").Value = "  +4.50%  "
